# Update of all values to match PDF edition 10 (commit 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (text unchanged, but keep explicit for clarity / safety)
$ws.Range("A1").Value = "Aircraft Type"
$ws.Range("B1").Value = "Flights"
$ws.Range("C1").Value = "Proportion"
$ws.Range("D1").Value = "Cumulative"

# Data rows, re-sorted/updated per new PDF edition figures
$data = @(
    @("B738", 1717381, 0.2,  0.2),
    @("A320", 1417930, 0.16, 0.36),
    @("A319", 513699,  0.06, 0.42),
    @("A20N", 462294,  0.05, 0.47),
    @("A321", 371456,  0.04, 0.51),
    @("A21N", 296722,  0.03, 0.55),
    @("B38M", 270614,  0.03, 0.58),
    @("AT76", 210778,  0.02, 0.6),
    @("E190", 210022,  0.02, 0.63),
    @("B77W", 163600,  0.02, 0.64),
    @("B789", 135989,  0.02, 0.66),
    @("A333", 128892,  0.01, 0.67),
    @("E195", 124823,  0.01, 0.69),
    @("CRJ9", 112652,  0.01, 0.7),
    @("AT75", 100619,  0.01, 0.71),
    @("DH8D", 93110,   0.01, 0.72),
    @("BCS3", 90613,   0.01, 0.73),
    @("A332", 86165,   0.01, 0.74),
    @("B788", 80654,   0.01, 0.75),
    @("A359", 78200,   0.01, 0.76),
    @("Other types", 2073792, 0.24, 1)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $row = $row + 1
}

# Total row: label unchanged, Flights becomes a plain value (formula removed), Proportion unchanged at 1
$ws.Range("A23").Value = "Total"
$ws.Range("B23").Value = 8740005
$ws.Range("C23").Value = 1

# Selection moved as part of the edit
[void]$ws.Range("I25").Select()
